$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, styled like the other header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# time_taken values for data rows 2-94
$timeTaken = @(
    "2021-10-05 10:52:45.497822",
    "2021-10-05 10:52:45.497833",
    "2021-10-05 10:52:45.497837",
    "2021-10-05 10:52:45.497839",
    "2021-10-05 10:52:45.497843",
    "2021-10-05 10:52:45.497845",
    "2021-10-05 10:52:45.497848",
    "2021-10-05 10:52:45.497850",
    "2021-10-05 10:52:45.497854",
    "2021-10-05 10:52:45.497856",
    "2021-10-05 10:52:45.497859",
    "2021-10-05 10:52:45.497861",
    "2021-10-05 10:52:45.497864",
    "2021-10-05 10:52:45.497866",
    "2021-10-05 10:52:45.497869",
    "2021-10-05 10:52:45.497872",
    "2021-10-05 10:52:45.497875",
    "2021-10-05 10:52:45.497877",
    "2021-10-05 10:52:45.497880",
    "2021-10-05 10:52:45.497883",
    "2021-10-05 10:52:45.497885",
    "2021-10-05 10:52:45.497888",
    "2021-10-05 10:52:45.497890",
    "2021-10-05 10:52:45.497893",
    "2021-10-05 10:52:45.497896",
    "2021-10-05 10:52:45.497898",
    "2021-10-05 10:52:45.497901",
    "2021-10-05 10:52:45.497904",
    "2021-10-05 10:52:45.497906",
    "2021-10-05 10:52:45.497909",
    "2021-10-05 10:52:45.497911",
    "2021-10-05 10:52:45.497914",
    "2021-10-05 10:52:45.497917",
    "2021-10-05 10:52:45.497920",
    "2021-10-05 10:52:45.497923",
    "2021-10-05 10:52:45.497925",
    "2021-10-05 10:52:45.497928",
    "2021-10-05 10:52:45.497930",
    "2021-10-05 10:52:45.497933",
    "2021-10-05 10:52:45.497936",
    "2021-10-05 10:52:45.497938",
    "2021-10-05 10:52:45.497941",
    "2021-10-05 10:52:45.497944",
    "2021-10-05 10:52:45.497947",
    "2021-10-05 10:52:45.497949",
    "2021-10-05 10:52:45.497952",
    "2021-10-05 10:52:45.497954",
    "2021-10-05 10:52:45.497957",
    "2021-10-05 10:52:45.497959",
    "2021-10-05 10:52:45.497962",
    "2021-10-05 10:52:45.497965",
    "2021-10-05 10:52:45.497967",
    "2021-10-05 10:52:45.497970",
    "2021-10-05 10:52:45.497973",
    "2021-10-05 10:52:45.497976",
    "2021-10-05 10:52:45.497978",
    "2021-10-05 10:52:45.497981",
    "2021-10-05 10:52:45.497983",
    "2021-10-05 10:52:45.497986",
    "2021-10-05 10:52:45.497989",
    "2021-10-05 10:52:45.497991",
    "2021-10-05 10:52:45.497994",
    "2021-10-05 10:52:45.497997",
    "2021-10-05 10:52:45.497999",
    "2021-10-05 10:52:45.498003",
    "2021-10-05 10:52:45.498006",
    "2021-10-05 10:52:45.498009",
    "2021-10-05 10:52:45.498012",
    "2021-10-05 10:52:45.498015",
    "2021-10-05 10:52:45.498017",
    "2021-10-05 10:52:45.498020",
    "2021-10-05 10:52:45.498023",
    "2021-10-05 10:52:45.498026",
    "2021-10-05 10:52:45.498029",
    "2021-10-05 10:52:45.498031",
    "2021-10-05 10:52:45.498034",
    "2021-10-05 10:52:45.498039",
    "2021-10-05 10:52:45.498042",
    "2021-10-05 10:52:45.498045",
    "2021-10-05 10:52:45.498048",
    "2021-10-05 10:52:45.498050",
    "2021-10-05 10:52:45.498054",
    "2021-10-05 10:52:45.498056",
    "2021-10-05 10:52:45.498059",
    "2021-10-05 10:52:45.498061",
    "2021-10-05 10:52:45.498064",
    "2021-10-05 10:52:45.498066",
    "2021-10-05 10:52:45.498069",
    "2021-10-05 10:52:45.498072",
    "2021-10-05 10:52:45.498075",
    "2021-10-05 10:52:45.498077",
    "2021-10-05 10:52:45.498080",
    "2021-10-05 10:52:45.498084"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

